# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E3) and "Correspond Handback DateTime" (H3)
# values for the a1bbab1d-... row on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-17 06:24:34"
$wsZhCn.Range("H3").Value = "2016-03-17 06:25:16"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-17 06:24:42"
$wsDeDe.Range("H3").Value = "2016-03-17 06:25:29"
